$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The table in A1:B10 lists "Real Madrid CF - <Opponent>" matches with dates.
# The first match (Valencia CF, 08/01/2022) is removed, and all the rows shift
# up by one. A new match (Real Betis Balompié, 22/05/2022) is appended at the
# end (row 10).

$matches = @(
    @("Real Madrid CF - Elche CF", "23/01/2022"),
    @("Real Madrid CF - Granada CF", "06/02/2022"),
    @("Real Madrid CF - Deportivo Alavés", "20/02/2022"),
    @("Real Madrid CF - Real Sociedad", "06/03/2022"),
    @("Real Madrid CF - Paris Saint-Germain", "09/03/2022"),
    @("Real Madrid CF - FC Barcelona", "20/03/2022"),
    @("Real Madrid CF - Getafe CF", "10/04/2022"),
    @("Real Madrid CF - RCD Espanyol", "01/05/2022"),
    @("Real Madrid CF - UD Levante", "11/05/2022"),
    @("Real Madrid CF - Real Betis Balompié", "22/05/2022")
)

# Simply assigning a plain string like "06/02/2022" via .Value/.Value2 gets
# auto-recognized by Excel as a date and silently converted into a date
# serial number with a new number-format style - which we do not want here
# (the target only changes the text content, not any formatting).
#
# To force the content to be stored as literal text (landing back in
# xl/sharedStrings.xml, unformatted) we enter it as a quoted-string formula
# (which always evaluates to text) and then replace the formula with its
# static result via copy / paste-special-values. This avoids touching
# xl/styles.xml entirely.
for ($i = 0; $i -lt $matches.Length; $i++) {
    $row = $i + 1
    $escapedName = $matches[$i][0] -replace '"', '""'
    $escapedDate = $matches[$i][1] -replace '"', '""'
    $ws.Cells.Item($row, 1).Formula = '="' + $escapedName + '"'
    $ws.Cells.Item($row, 2).Formula = '="' + $escapedDate + '"'
}

$used = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item($matches.Length, 2))
$used.Copy() | Out-Null
$used.PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = 0
